$d = $word.ActiveDocument

$d.Content.Find.Execute("859×7=6013", $true, $false, $false, $false, $false, $true, 1, $false, "325×9=2925", 2) | Out-Null
$d.Content.Find.Execute("528×5=2640", $true, $false, $false, $false, $false, $true, 1, $false, "208×8=1664", 2) | Out-Null
$d.Content.Find.Execute("551×5=2755", $true, $false, $false, $false, $false, $true, 1, $false, "867×6=5202", 2) | Out-Null
$d.Content.Find.Execute("458×2=916", $true, $false, $false, $false, $false, $true, 1, $false, "967×4=3868", 2) | Out-Null
$d.Content.Find.Execute("818×5=4090", $true, $false, $false, $false, $false, $true, 1, $false, "545×2=1090", 2) | Out-Null
$d.Content.Find.Execute("625×3=1875", $true, $false, $false, $false, $false, $true, 1, $false, "158×6=948", 2) | Out-Null
$d.Content.Find.Execute("476×3=1428", $true, $false, $false, $false, $false, $true, 1, $false, "726×7=5082", 2) | Out-Null
$d.Content.Find.Execute("345×3=1035", $true, $false, $false, $false, $false, $true, 1, $false, "276×4=1104", 2) | Out-Null
$d.Content.Find.Execute("827×8=6616", $true, $false, $false, $false, $false, $true, 1, $false, "392×4=1568", 2) | Out-Null
$d.Content.Find.Execute("872×3=2616", $true, $false, $false, $false, $false, $true, 1, $false, "268×5=1340", 2) | Out-Null
$d.Content.Find.Execute("812×3=2436", $true, $false, $false, $false, $false, $true, 1, $false, "677×7=4739", 2) | Out-Null
$d.Content.Find.Execute("365×8=2920", $true, $false, $false, $false, $false, $true, 1, $false, "602×5=3010", 2) | Out-Null
$d.Content.Find.Execute("885×7=6195", $true, $false, $false, $false, $false, $true, 1, $false, "142×5=710", 2) | Out-Null
$d.Content.Find.Execute("481×8=3848", $true, $false, $false, $false, $false, $true, 1, $false, "652×2=1304", 2) | Out-Null
$d.Content.Find.Execute("683×7=4781", $true, $false, $false, $false, $false, $true, 1, $false, "815×4=3260", 2) | Out-Null
$d.Content.Find.Execute("102×5=510", $true, $false, $false, $false, $false, $true, 1, $false, "517×4=2068", 2) | Out-Null
$d.Content.Find.Execute("659×2=1318", $true, $false, $false, $false, $false, $true, 1, $false, "535×2=1070", 2) | Out-Null
$d.Content.Find.Execute("988×9=8892", $true, $false, $false, $false, $false, $true, 1, $false, "494×3=1482", 2) | Out-Null
$d.Content.Find.Execute("624×3=1872", $true, $false, $false, $false, $false, $true, 1, $false, "143×8=1144", 2) | Out-Null
$d.Content.Find.Execute("665×3=1995", $true, $false, $false, $false, $false, $true, 1, $false, "750×6=4500", 2) | Out-Null
$d.Content.Find.Execute("154×9=1386", $true, $false, $false, $false, $false, $true, 1, $false, "346×2=692", 2) | Out-Null
$d.Content.Find.Execute("838×4=3352", $true, $false, $false, $false, $false, $true, 1, $false, "536×3=1608", 2) | Out-Null
$d.Content.Find.Execute("365×9=3285", $true, $false, $false, $false, $false, $true, 1, $false, "380×4=1520", 2) | Out-Null
$d.Content.Find.Execute("465×2=930", $true, $false, $false, $false, $false, $true, 1, $false, "345×7=2415", 2) | Out-Null
$d.Content.Find.Execute("470×7=3290", $true, $false, $false, $false, $false, $true, 1, $false, "516×3=1548", 2) | Out-Null
